$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'273.36"
$ws.Range("E2").Value = "'0.54%"
$ws.Range("D3").Value = "'26.80"
$ws.Range("E3").Value = "'0.20%"
$ws.Range("D4").Value = "'4.901"
$ws.Range("E4").Value = "'3.89%"
$ws.Range("D5").Value = "'0.06321"
$ws.Range("E5").Value = "'3.04%"
$ws.Range("D6").Value = "'6.898"
$ws.Range("E6").Value = "'2.33%"
$ws.Range("D7").Value = "'3.359"
$ws.Range("E7").Value = "'5.57%"
$ws.Range("D8").Value = "'1.342"
$ws.Range("E8").Value = "'48.68%"
$ws.Range("D9").Value = "'0.8840"
$ws.Range("E9").Value = "'3.34%"
$ws.Range("D10").Value = "'0.1465"
$ws.Range("E10").Value = "'2.22%"
$ws.Range("D11").Value = "'0.05097"
$ws.Range("E11").Value = "'1.23%"
$ws.Range("D12").Value = "'0.07403"
$ws.Range("E12").Value = "'3.62%"
$ws.Range("D13").Value = "'0.03145"
$ws.Range("E13").Value = "'-0.62%"
$ws.Range("D14").Value = "'0.09036"
$ws.Range("E14").Value = "'0.00%"
$ws.Range("D15").Value = "'0.001570"
$ws.Range("E15").Value = "'2.88%"
$ws.Range("D16").Value = "'0.0006295"
$ws.Range("E16").Value = "'3.92%"
$ws.Range("D17").Value = "'0.006025"
$ws.Range("E17").Value = "'-0.22%"
$ws.Range("D18").Value = "'3.461"
$ws.Range("E18").Value = "'-0.02%"
$ws.Range("E19").Value = "'-0.25%"
$ws.Range("E20").Value = "'2.63%"
$ws.Range("D21").Value = "'0.1334"
$ws.Range("E21").Value = "'4.07%"
$ws.Range("D22").Value = "'3.907"
$ws.Range("E22").Value = "'1.54%"
$ws.Range("D23").Value = "'0.04332"
$ws.Range("E23").Value = "'2.05%"
$ws.Range("D24").Value = "'0.001178"
$ws.Range("E24").Value = "'-0.12%"
$ws.Range("E25").Value = "'-12.04%"
$ws.Range("D26").Value = "'0.0001201"
$ws.Range("E26").Value = "'0.06%"
$ws.Range("D27").Value = "'0.0001699"
$ws.Range("E27").Value = "'1.31%"
$ws.Range("D40").Value = "'0.04036"
$ws.Range("E40").Value = "'1.60%"
$ws.Range("D41").Value = "'0.006613"
$ws.Range("E41").Value = "'57.51%"
$ws.Range("D42").Value = "'0.1164"
$ws.Range("E42").Value = "'4.08%"
$ws.Range("D43").Value = "'0.002133"
$ws.Range("E43").Value = "'3.50%"
$ws.Range("D44").Value = "'0.01227"
$ws.Range("E44").Value = "'5.09%"
$ws.Range("D45").Value = "'0.00005349"
$ws.Range("E45").Value = "'3.86%"
$ws.Range("E46").Value = "'162.04%"
$ws.Range("D47").Value = "'0.02123"
$ws.Range("E47").Value = "'-29.07%"
